$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Cells.Item(1,1).Value = "Name"
$ws.Cells.Item(1,2).Value = "Address"

# --- Patient rows (alphabetical by first name) ---
$ws.Cells.Item(2,1).Value = "Alice Person"
$ws.Cells.Item(2,2).Value = "480 Alta Rd, San Diego CA 92179"

$ws.Cells.Item(3,1).Value = "Barry Barber"
$ws.Cells.Item(3,2).Value = "220 Euclid Ave #120, San Diego, CA 92114"

$ws.Cells.Item(4,1).Value = "Carl Person"
$ws.Cells.Item(4,2).Value = "1173 Front St, San Diego, CA 92101"

$ws.Cells.Item(5,1).Value = "Diego Person"
$ws.Cells.Item(5,2).Value = "201 N Broadway, Escondido CA 92025"

$ws.Cells.Item(6,1).Value = "Edgar Person"
$ws.Cells.Item(6,2).Value = "233 4th Ave, Chula Vista, CA 91910"

$ws.Cells.Item(7,1).Value = "Frances Person"
$ws.Cells.Item(7,2).Value = "3232 Main St, Lemon Grove, CA 91945"

$ws.Cells.Item(8,1).Value = "Gail Person"
$ws.Cells.Item(8,2).Value = "1243 National City Blvd, National City, CA 91950"

$ws.Cells.Item(9,1).Value = "Harriet Person"
$ws.Cells.Item(9,2).Value = "400 South Melrose Drive, Suite 108 Vista, CA 92081"

# --- Header formatting: bold, centered (built on a scratch cell so the
#     final header style lands in a single cellXf, mirroring how the
#     workbook was authored) ---
$scratch = $ws.Range("Z1")
$scratch.Font.Bold = $true
$scratch.HorizontalAlignment = -4108
$scratch.VerticalAlignment = -4108
$scratch.Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)
$scratch.Clear()

# --- Column widths ---
$ws.Columns("A:A").ColumnWidth = 23.2
$ws.Columns("B:E").ColumnWidth = 46.3

# --- Page setup ---
$ws.PageSetup.Orientation = 1

# --- Selection ---
$ws.Range("B16").Select()
